$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting of column A (bold, bordered, centered style) down to the newly added rows 24-43
# by copying the format already used in A2:A23.
$ws.Range("A23").Copy()
$ws.Range("A24:A43").PasteSpecial(-4122)

# Column A holds the running index (row - 2) for every data row, rows 2 through 43
for ($r = 2; $r -le 43; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# Update the stock ticker columns (B..F) that changed between the old and new lists
$ws.Range("B2").Value = "NSE:AVALON"
$ws.Range("C2").Value = "NSE:ADANIPORTS"
$ws.Range("D2").Value = "NSE:ATUL"
$ws.Range("E2").Value = "NSE:ABB"
$ws.Range("F2").Value = "NSE:COFORGE"
$ws.Range("B3").Value = "NSE:BIRLAMONEY"
$ws.Range("C3").Value = "NSE:AMJLAND"
$ws.Range("E3").Value = "NSE:ABCAPITAL"
$ws.Range("F3").Value = ""
$ws.Range("B4").Value = "NSE:COFORGE"
$ws.Range("C4").Value = "NSE:ANANDRATHI"
$ws.Range("E4").Value = "NSE:ABFRL"
$ws.Range("F4").Value = ""
$ws.Range("B5").Value = "NSE:SAGARDEEP"
$ws.Range("C5").Value = "NSE:APEX"
$ws.Range("E5").Value = "NSE:BANKBARODA"
$ws.Range("F5").Value = ""
$ws.Range("B6").Value = ""
$ws.Range("C6").Value = "NSE:APOLLOTYRE"
$ws.Range("E6").Value = "NSE:BANKINDIA"
$ws.Range("B7").Value = ""
$ws.Range("C7").Value = "NSE:APTECHT"
$ws.Range("E7").Value = "NSE:BATAINDIA"
$ws.Range("B8").Value = ""
$ws.Range("C8").Value = "NSE:ASTRAL"
$ws.Range("E8").Value = "NSE:BHARATFORG"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "NSE:AVANTIFEED"
$ws.Range("E9").Value = "NSE:DMART"
$ws.Range("B10").Value = ""
$ws.Range("C10").Value = "NSE:BEML"
$ws.Range("E10").Value = "NSE:HAL"
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = "NSE:CELEBRITY"
$ws.Range("E11").Value = "NSE:ICICIBANK"
$ws.Range("B12").Value = ""
$ws.Range("C12").Value = "NSE:CREATIVE"
$ws.Range("E12").Value = "NSE:IRCTC"
$ws.Range("B13").Value = ""
$ws.Range("C13").Value = "NSE:DELHIVERY"
$ws.Range("E13").Value = "NSE:IRFC"
$ws.Range("B14").Value = ""
$ws.Range("C14").Value = "NSE:DLF"
$ws.Range("E14").Value = "NSE:JSWSTEEL"
$ws.Range("B15").Value = ""
$ws.Range("C15").Value = "NSE:DPWIRES"
$ws.Range("E15").Value = "NSE:LTIM"
$ws.Range("B16").Value = ""
$ws.Range("C16").Value = "NSE:DREDGECORP"
$ws.Range("E16").Value = "NSE:MARUTI"
$ws.Range("B17").Value = ""
$ws.Range("C17").Value = "NSE:EPL"
$ws.Range("E17").Value = "NSE:MRF"
$ws.Range("B18").Value = ""
$ws.Range("C18").Value = "NSE:EXXARO"
$ws.Range("E18").Value = "NSE:NTPC"
$ws.Range("B19").Value = ""
$ws.Range("C19").Value = "NSE:GNFC"
$ws.Range("E19").Value = "NSE:PEL"
$ws.Range("B20").Value = ""
$ws.Range("C20").Value = "NSE:HBLPOWER"
$ws.Range("E20").Value = "NSE:PIDILITIND"
$ws.Range("C21").Value = "NSE:HFCL"
$ws.Range("C22").Value = "NSE:HINDWAREAP"
$ws.Range("C23").Value = "NSE:HUDCO"
$ws.Range("C24").Value = "NSE:IFGLEXPOR"
$ws.Range("C25").Value = "NSE:INFOBEAN"
$ws.Range("C26").Value = "NSE:INOXWIND"
$ws.Range("C27").Value = "NSE:MAHLOG"
$ws.Range("C28").Value = "NSE:MANINDS"
$ws.Range("C29").Value = "NSE:METROPOLIS"
$ws.Range("C30").Value = "NSE:MIRZAINT"
$ws.Range("C31").Value = "NSE:MRF"
$ws.Range("C32").Value = "NSE:NAHARINDUS"
$ws.Range("C33").Value = "NSE:NBCC"
$ws.Range("C34").Value = "NSE:NIITLTD"
$ws.Range("C35").Value = "NSE:ORIENTCER"
$ws.Range("C36").Value = "NSE:PARAGMILK"
$ws.Range("C37").Value = "NSE:PATINTLOG"
$ws.Range("C38").Value = "NSE:PDSL"
$ws.Range("C39").Value = "NSE:PNBGILTS"
$ws.Range("C40").Value = "NSE:RAMCOCEM"
$ws.Range("C41").Value = "NSE:ROHLTD"
$ws.Range("C42").Value = "NSE:ROTO"
$ws.Range("C43").Value = "NSE:RUPA"
